# Build the target values for D, J, K, L, M, P per row (rows 2-28, row 8 unchanged)
$rowValues = @{
    2 = @{ D = 44463; J = 25; K = 24000; L = 25000; M = 24480; P = 1632 }
    3 = @{ D = 44425; J = 25; K = 24000; L = 25000; M = 24520; P = 1635 }
    4 = @{ D = 44428; J = 16; K = 25000; L = 26000; M = 25500; P = 1700 }
    5 = @{ D = 44432; J = 34; K = 24000; L = 25000; M = 24500; P = 1633 }
    6 = @{ D = 44453; J = 25; K = 25000; L = 26000; M = 25520; P = 1701 }
    7 = @{ D = 44449; J = 18; K = 24000; L = 25000; M = 24500; P = 1633 }
    9 = @{ D = 44435; J = 34; K = 24000; L = 25000; M = 24500; P = 1633 }
    10 = @{ D = 44406; J = 25; K = 24000; L = 25000; M = 24520; P = 1635 }
    11 = @{ D = 44341; J = 36; K = 24000; L = 25000; M = 24500; P = 1633 }
    12 = @{ D = 44442; J = 28; K = 24000; L = 25000; M = 24500; P = 1633 }
    13 = @{ D = 44460; J = 25; K = 24000; L = 25000; M = 24480; P = 1632 }
    14 = @{ D = 44351; J = 34; K = 24000; L = 25000; M = 24500; P = 1633 }
    15 = @{ D = 44411; J = 34; K = 25000; L = 26000; M = 25500; P = 1700 }
    16 = @{ D = 44343; J = 26; K = 23000; L = 24000; M = 23500; P = 1567 }
    17 = @{ D = 44336; J = 34; K = 24000; L = 25000; M = 24500; P = 1633 }
    18 = @{ D = 44400; J = 16; K = 24000; L = 25000; M = 24500; P = 1633 }
    19 = @{ D = 44455; J = 18; K = 24000; L = 25000; M = 24500; P = 1633 }
    20 = @{ D = 44397; J = 34; K = 23000; L = 24000; M = 23500; P = 1567 }
    21 = @{ D = 44446; J = 34; K = 24000; L = 25000; M = 24500; P = 1633 }
    22 = @{ D = 44421; J = 18; K = 24000; L = 25000; M = 24500; P = 1633 }
    23 = @{ D = 44329; J = 25; K = 23000; L = 23000; M = 23000; P = 1533 }
    24 = @{ D = 44385; J = 25; K = 14000; L = 15000; M = 14480; P = 965 }
    25 = @{ D = 44413; J = 25; K = 24000; L = 25000; M = 24480; P = 1632 }
    26 = @{ D = 44390; J = 34; K = 24000; L = 25000; M = 24500; P = 1633 }
    27 = @{ D = 44349; J = 21; K = 24000; L = 25000; M = 24524; P = 1635 }
    28 = @{ D = 44418; J = 16; K = 25000; L = 26000; M = 25500; P = 1700 }
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($row in $rowValues.Keys) {
    $vals = $rowValues[$row]
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("J$row").Value = $vals.J
    $ws.Range("K$row").Value = $vals.K
    $ws.Range("L$row").Value = $vals.L
    $ws.Range("M$row").Value = $vals.M
    $ws.Range("P$row").Value = $vals.P
}
